# chore: update Sheets via scheduled runner
# Refresh the currentAveragePrice / LevePrice / LeveProfit market-data
# columns (H:N) for the affected leve rows across each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4176.4346
$ws.Range("I40").Value = 3282.75
$ws.Range("J40").Value = 5151.364
$ws.Range("K40").Value = 3282.75
$ws.Range("L40").Value = 5151.364
$ws.Range("M40").Value = -3107.75
$ws.Range("N40").Value = -5501.364

$ws.Range("H86").Value = 1775
$ws.Range("J86").Value = 1835.5
$ws.Range("L86").Value = 1835.5
$ws.Range("N86").Value = -4081.5

$ws.Range("H89").Value = 1775
$ws.Range("J89").Value = 1835.5
$ws.Range("L89").Value = 9177.5
$ws.Range("N89").Value = -20409.5

$ws.Range("H100").Value = 5599.4
$ws.Range("I100").Value = 3999
$ws.Range("K100").Value = 3999
$ws.Range("M100").Value = -3458

$ws.Range("H113").Value = 1476.3334
$ws.Range("J113").Value = 1324
$ws.Range("L113").Value = 1324
$ws.Range("N113").Value = -7832

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3393.25
$ws.Range("I45").Value = 2778.2222
$ws.Range("J45").Value = 4184
$ws.Range("K45").Value = 2778.2222
$ws.Range("L45").Value = 4184
$ws.Range("M45").Value = -2401.2222
$ws.Range("N45").Value = -4938

$ws.Range("H122").Value = 1241.6666
$ws.Range("J122").Value = 1241.6666
$ws.Range("L122").Value = 3724.9998
$ws.Range("N122").Value = -8624.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3929
$ws.Range("J20").Value = 1988
$ws.Range("L20").Value = 1988
$ws.Range("N20").Value = -2482

$ws.Range("H22").Value = 334
$ws.Range("I22").Value = 350
$ws.Range("J22").Value = 310
$ws.Range("K22").Value = 350
$ws.Range("L22").Value = 310
$ws.Range("M22").Value = -177
$ws.Range("N22").Value = -656

$ws.Range("H82").Value = 38779.777
$ws.Range("I82").Value = 19656.2
$ws.Range("J82").Value = 62684.25
$ws.Range("K82").Value = 19656.2
$ws.Range("L82").Value = 62684.25
$ws.Range("M82").Value = -19273.2
$ws.Range("N82").Value = -63450.25

$ws.Range("H85").Value = 38779.777
$ws.Range("I85").Value = 19656.2
$ws.Range("J85").Value = 62684.25
$ws.Range("K85").Value = 19656.2
$ws.Range("L85").Value = 62684.25
$ws.Range("M85").Value = -18330.2
$ws.Range("N85").Value = -65336.25

$ws.Range("H135").Value = 195000
$ws.Range("J135").Value = 195000
$ws.Range("L135").Value = 195000
$ws.Range("N135").Value = -205140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 6293.5
$ws.Range("I15").Value = 175
$ws.Range("J15").Value = 8333
$ws.Range("K15").Value = 175
$ws.Range("L15").Value = 8333
$ws.Range("M15").Value = -5
$ws.Range("N15").Value = -8673

$ws.Range("H31").Value = 5776.263
$ws.Range("I31").Value = 3914.5
$ws.Range("J31").Value = 8967.857
$ws.Range("K31").Value = 3914.5
$ws.Range("L31").Value = 8967.857
$ws.Range("M31").Value = -3619.5
$ws.Range("N31").Value = -9557.857

$ws.Range("H34").Value = 5776.263
$ws.Range("I34").Value = 3914.5
$ws.Range("J34").Value = 8967.857
$ws.Range("K34").Value = 3914.5
$ws.Range("L34").Value = 8967.857
$ws.Range("M34").Value = -3712.5
$ws.Range("N34").Value = -9371.857

$ws.Range("H122").Value = 2563.4443
$ws.Range("I122").Value = 2324
$ws.Range("J122").Value = 2862.75
$ws.Range("K122").Value = 6972
$ws.Range("L122").Value = 8588.25
$ws.Range("M122").Value = -4522
$ws.Range("N122").Value = -13488.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5101.4346
$ws.Range("J80").Value = 5388.8335
$ws.Range("L80").Value = 16166.5005
$ws.Range("N80").Value = -18038.5005

$ws.Range("H83").Value = 5101.4346
$ws.Range("J83").Value = 5388.8335
$ws.Range("L83").Value = 48499.5015
$ws.Range("N83").Value = -57859.5015

$ws.Range("H131").Value = 1210
$ws.Range("J131").Value = 1420
$ws.Range("L131").Value = 4260
$ws.Range("N131").Value = -14340

$ws.Range("H137").Value = 2874.5
$ws.Range("J137").Value = 3000
$ws.Range("L137").Value = 9000
$ws.Range("N137").Value = -19200

$ws.Range("H141").Value = 3432.5
$ws.Range("I141").Value = 1910
$ws.Range("J141").Value = 8000
$ws.Range("K141").Value = 5730
$ws.Range("L141").Value = 24000
$ws.Range("M141").Value = -550
$ws.Range("N141").Value = -34360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 138.2
$ws.Range("I2").Value = 16
$ws.Range("J2").Value = 168.75
$ws.Range("K2").Value = 16
$ws.Range("L2").Value = 168.75
$ws.Range("M2").Value = 97
$ws.Range("N2").Value = -394.75

$ws.Range("H70").Value = 9082.333000000001
$ws.Range("I70").Value = 8499.666999999999
$ws.Range("J70").Value = 9665
$ws.Range("K70").Value = 8499.666999999999
$ws.Range("L70").Value = 9665
$ws.Range("M70").Value = -8229.666999999999
$ws.Range("N70").Value = -10205

$ws.Range("H73").Value = 9082.333000000001
$ws.Range("I73").Value = 8499.666999999999
$ws.Range("J73").Value = 9665
$ws.Range("K73").Value = 8499.666999999999
$ws.Range("L73").Value = 9665
$ws.Range("M73").Value = -7563.666999999999
$ws.Range("N73").Value = -11537

$ws.Range("H122").Value = 2519.5557
$ws.Range("I122").Value = 1862.375
$ws.Range("J122").Value = 7777
$ws.Range("K122").Value = 5587.125
$ws.Range("L122").Value = 23331
$ws.Range("M122").Value = -3137.125
$ws.Range("N122").Value = -28231

$ws.Range("H132").Value = 123445.78
$ws.Range("I132").Value = 177670.17
$ws.Range("J132").Value = 14997
$ws.Range("K132").Value = 533010.51
$ws.Range("L132").Value = 44991
$ws.Range("M132").Value = -530480.51
$ws.Range("N132").Value = -50051

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1005.36365
$ws.Range("I22").Value = 642
$ws.Range("K22").Value = 642
$ws.Range("M22").Value = -347

$ws.Range("H27").Value = 1005.36365
$ws.Range("I27").Value = 642
$ws.Range("K27").Value = 642
$ws.Range("M27").Value = -535

$ws.Range("H34").Value = 3334933.2
$ws.Range("I34").Value = 3334933.2
$ws.Range("K34").Value = 3334933.2
$ws.Range("M34").Value = -3334761.2

$ws.Range("H40").Value = 8466.267
$ws.Range("I40").Value = 7908.636
$ws.Range("K40").Value = 7908.636
$ws.Range("M40").Value = -7772.636

$ws.Range("H122").Value = 5123.75
$ws.Range("I122").Value = 3998.3333
$ws.Range("K122").Value = 11994.9999
$ws.Range("M122").Value = -9544.999899999999

$ws.Range("H132").Value = 5998.5
$ws.Range("I132").Value = 4398.4
$ws.Range("J132").Value = 8665.333000000001
$ws.Range("K132").Value = 13195.2
$ws.Range("L132").Value = 25995.999
$ws.Range("M132").Value = -10665.2
$ws.Range("N132").Value = -31055.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3173.3333
$ws.Range("I122").Value = 2208.5881
$ws.Range("K122").Value = 6625.7643
$ws.Range("M122").Value = -4175.7643
